$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H8").Value = "testing"
